# Update the Manchester City fixture list:
# - Remove the two already-played fixtures (Swindon Town FA Cup replay, and
#   the Chelsea Premier League match) that used to occupy rows 1-2.
# - Keep the Brentford fixture that was previously in row 3, moving it to row 1.
# - Insert two newly scheduled fixtures (Norwich City away, then Sporting CP
#   in the Champions League) in rows 2-3.
# - Leave the remaining fixtures (rows 4 onward, i.e. what used to be rows
#   4-10) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Manchester City v Brentford"
$ws.Range("B1").Value = "09/02/2022 19:45 | Premier League"

$ws.Range("A2").Value = "Norwich City v Manchester City"
$ws.Range("B2").Value = "12/02/2022 17:30 | Premier League"

$ws.Range("A3").Value = "Sporting CP v Manchester City"
$ws.Range("B3").Value = "15/02/2022 20:00 | UEFA Champions League"
